$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A69").Value = 0
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 0.102041
$ws.Range("E69").Value = -23.10956470571409
$ws.Range("F69").Value = "query"

$ws.Range("A70").Value = 0
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 0.204081
$ws.Range("E70").Value = -21.58417940885292
$ws.Range("F70").Value = "query"
